$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3891.96
$ws.Range("I64").Value = 3376.4707
$ws.Range("J64").Value = 4987.375
$ws.Range("K64").Value = 3376.4707
$ws.Range("L64").Value = 4987.375
$ws.Range("M64").Value = -3128.4707
$ws.Range("N64").Value = -5483.375
$ws.Range("H67").Value = 3891.96
$ws.Range("I67").Value = 3376.4707
$ws.Range("J67").Value = 4987.375
$ws.Range("K67").Value = 3376.4707
$ws.Range("L67").Value = 4987.375
$ws.Range("M67").Value = -2518.4707
$ws.Range("N67").Value = -6703.375
$ws.Range("H74").Value = 5162.273
$ws.Range("I74").Value = 4455.143
$ws.Range("J74").Value = 6399.75
$ws.Range("K74").Value = 4455.143
$ws.Range("L74").Value = 6399.75
$ws.Range("M74").Value = -3519.143
$ws.Range("N74").Value = -8271.75
$ws.Range("H77").Value = 5162.273
$ws.Range("I77").Value = 4455.143
$ws.Range("J77").Value = 6399.75
$ws.Range("K77").Value = 22275.715
$ws.Range("L77").Value = 31998.75
$ws.Range("M77").Value = -17595.715
$ws.Range("N77").Value = -41358.75
$ws.Range("H116").Value = 2032.5
$ws.Range("I116").Value = 1739.8334
$ws.Range("J116").Value = 2325.1667
$ws.Range("K116").Value = 1739.8334
$ws.Range("L116").Value = 2325.1667
$ws.Range("M116").Value = 1702.1666
$ws.Range("N116").Value = -9209.1667
$ws.Range("H132").Value = 1847.8387
$ws.Range("I132").Value = 1722.4231
$ws.Range("K132").Value = 5167.2693
$ws.Range("M132").Value = -2637.2693
$ws.Range("H135").Value = 51725216
$ws.Range("I135").Value = 20834476
$ws.Range("J135").Value = 200000770
$ws.Range("K135").Value = 187510284
$ws.Range("L135").Value = 1800006930
$ws.Range("M135").Value = -187507749
$ws.Range("N135").Value = -1800012000
$ws.Range("H136").Value = 63111.668
$ws.Range("J136").Value = 63111.668
$ws.Range("L136").Value = 63111.668
$ws.Range("N136").Value = -73311.66800000001
$ws.Range("H137").Value = 2632.3872
$ws.Range("I137").Value = 2643.5217
$ws.Range("J137").Value = 2600.375
$ws.Range("K137").Value = 7930.5651
$ws.Range("L137").Value = 7801.125
$ws.Range("M137").Value = -5380.5651
$ws.Range("N137").Value = -12901.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 42117.46
$ws.Range("I32").Value = 45464.75
$ws.Range("J32").Value = 1950
$ws.Range("K32").Value = 45464.75
$ws.Range("L32").Value = 1950
$ws.Range("M32").Value = -45177.75
$ws.Range("N32").Value = -2524
$ws.Range("H74").Value = 16333.625
$ws.Range("I74").Value = 4931
$ws.Range("J74").Value = 35338
$ws.Range("K74").Value = 4931
$ws.Range("L74").Value = 35338
$ws.Range("M74").Value = -4057
$ws.Range("N74").Value = -37086
$ws.Range("H77").Value = 16333.625
$ws.Range("I77").Value = 4931
$ws.Range("J77").Value = 35338
$ws.Range("K77").Value = 24655
$ws.Range("L77").Value = 176690
$ws.Range("M77").Value = -20287
$ws.Range("N77").Value = -185426
$ws.Range("H132").Value = 3170.913
$ws.Range("I132").Value = 2805.9333
$ws.Range("J132").Value = 3855.25
$ws.Range("K132").Value = 8417.7999
$ws.Range("L132").Value = 11565.75
$ws.Range("M132").Value = -5887.7999
$ws.Range("N132").Value = -16625.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2843.3333
$ws.Range("I16").Value = 1815
$ws.Range("K16").Value = 1815
$ws.Range("M16").Value = -1528
$ws.Range("H31").Value = 6727.553
$ws.Range("I31").Value = 7151.353
$ws.Range("K31").Value = 7151.353
$ws.Range("M31").Value = -6856.353
$ws.Range("H34").Value = 6727.553
$ws.Range("I34").Value = 7151.353
$ws.Range("K34").Value = 7151.353
$ws.Range("M34").Value = -6949.353
$ws.Range("H113").Value = 2843.3333
$ws.Range("I113").Value = 1815
$ws.Range("K113").Value = 1815
$ws.Range("M113").Value = 355
$ws.Range("H122").Value = 9701.380999999999
$ws.Range("I122").Value = 5081.3125
$ws.Range("K122").Value = 15243.9375
$ws.Range("M122").Value = -12793.9375
$ws.Range("H132").Value = 5333.0586
$ws.Range("I132").Value = 6099.92
$ws.Range("J132").Value = 3202.889
$ws.Range("K132").Value = 18299.76
$ws.Range("L132").Value = 9608.667000000001
$ws.Range("M132").Value = -15769.76
$ws.Range("N132").Value = -14668.667
$ws.Range("H134").Value = 2710
$ws.Range("I134").Value = 2044.2354
$ws.Range("J134").Value = 4767.8184
$ws.Range("K134").Value = 6132.706200000001
$ws.Range("L134").Value = 14303.4552
$ws.Range("M134").Value = -3597.706200000001
$ws.Range("N134").Value = -19373.4552

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1879.0416
$ws.Range("I34").Value = 327.33334
$ws.Range("J34").Value = 2810.0667
$ws.Range("K34").Value = 982.0000200000001
$ws.Range("L34").Value = 8430.2001
$ws.Range("M34").Value = -898.0000200000001
$ws.Range("N34").Value = -8598.2001
$ws.Range("H39").Value = 5490.5264
$ws.Range("J39").Value = 5490.5264
$ws.Range("L39").Value = 16471.5792
$ws.Range("N39").Value = -17059.5792
$ws.Range("H55").Value = 3722.2222
$ws.Range("J55").Value = 3722.2222
$ws.Range("L55").Value = 11166.6666
$ws.Range("N55").Value = -11520.6666
$ws.Range("H113").Value = 678.25354
$ws.Range("I113").Value = 691.11865
$ws.Range("J113").Value = 615
$ws.Range("K113").Value = 2073.35595
$ws.Range("L113").Value = 1845
$ws.Range("M113").Value = 96.64404999999988
$ws.Range("N113").Value = -6185

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6926.3335
$ws.Range("I80").Value = 7849.5
$ws.Range("J80").Value = 5080
$ws.Range("K80").Value = 7849.5
$ws.Range("L80").Value = 5080
$ws.Range("M80").Value = -6851.5
$ws.Range("N80").Value = -7076
$ws.Range("H83").Value = 6926.3335
$ws.Range("I83").Value = 7849.5
$ws.Range("J83").Value = 5080
$ws.Range("K83").Value = 39247.5
$ws.Range("L83").Value = 25400
$ws.Range("M83").Value = -34255.5
$ws.Range("N83").Value = -35384
$ws.Range("H122").Value = 10535.143
$ws.Range("I122").Value = 15249
$ws.Range("J122").Value = 4250
$ws.Range("K122").Value = 45747
$ws.Range("L122").Value = 12750
$ws.Range("M122").Value = -43297
$ws.Range("N122").Value = -17650
$ws.Range("H132").Value = 11217.296
$ws.Range("I132").Value = 9483.883
$ws.Range("J132").Value = 14164.1
$ws.Range("K132").Value = 28451.649
$ws.Range("L132").Value = 42492.3
$ws.Range("M132").Value = -25921.649
$ws.Range("N132").Value = -47552.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 603.3333
$ws.Range("J22").Value = 375
$ws.Range("L22").Value = 375
$ws.Range("N22").Value = -965
$ws.Range("H27").Value = 603.3333
$ws.Range("J27").Value = 375
$ws.Range("L27").Value = 375
$ws.Range("N27").Value = -589
$ws.Range("H40").Value = 2956
$ws.Range("I40").Value = 2888.8235
$ws.Range("J40").Value = 3336.6667
$ws.Range("K40").Value = 2888.8235
$ws.Range("L40").Value = 3336.6667
$ws.Range("M40").Value = -2752.8235
$ws.Range("N40").Value = -3608.6667
$ws.Range("H46").Value = 1220
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 775
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 775
$ws.Range("M46").Value = -2812
$ws.Range("N46").Value = -1151
$ws.Range("H55").Value = 281.16666
$ws.Range("I55").Value = 258
$ws.Range("J55").Value = 297.7143
$ws.Range("K55").Value = 258
$ws.Range("L55").Value = 297.7143
$ws.Range("M55").Value = -85
$ws.Range("N55").Value = -643.7143
$ws.Range("H122").Value = 7759.3105
$ws.Range("I122").Value = 7175.75
$ws.Range("J122").Value = 9056.111000000001
$ws.Range("K122").Value = 21527.25
$ws.Range("L122").Value = 27168.333
$ws.Range("M122").Value = -19077.25
$ws.Range("N122").Value = -32068.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3534
$ws.Range("I62").Value = 3800
$ws.Range("J62").Value = 3002
$ws.Range("K62").Value = 3800
$ws.Range("L62").Value = 3002
$ws.Range("M62").Value = -3176
$ws.Range("N62").Value = -4250
$ws.Range("H65").Value = 3534
$ws.Range("I65").Value = 3800
$ws.Range("J65").Value = 3002
$ws.Range("K65").Value = 19000
$ws.Range("L65").Value = 15010
$ws.Range("M65").Value = -15880
$ws.Range("N65").Value = -21250
$ws.Range("H122").Value = 4284.086
$ws.Range("I122").Value = 1397.9231
$ws.Range("J122").Value = 12621.889
$ws.Range("K122").Value = 4193.7693
$ws.Range("L122").Value = 37865.667
$ws.Range("M122").Value = -1743.7693
$ws.Range("N122").Value = -42765.667
$ws.Range("H132").Value = 3908.0952
$ws.Range("I132").Value = 3439.1333
$ws.Range("J132").Value = 5080.5
$ws.Range("K132").Value = 10317.3999
$ws.Range("L132").Value = 15241.5
$ws.Range("M132").Value = -7787.3999
$ws.Range("N132").Value = -20301.5
$ws.Range("H136").Value = 5209
$ws.Range("J136").Value = 9192.137000000001
$ws.Range("L136").Value = 27576.411
$ws.Range("N136").Value = -32676.411
